$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '28.315.93'
$c.Style = $s
$ws.Range("E2").Value = '  -0.46%  '
$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.808.30'
$c.Style = $s
$ws.Range("E4").Value = '  -0.09%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '312.56'
$c.Style = $s
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("E6").Value = '  -0.13%  '
$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.5149'
$c.Style = $s
$ws.Range("E7").Value = '  -0.39%  '
$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.3988'
$c.Style = $s
$ws.Range("E8").Value = '  +3.55%  '
$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.07851'
$c.Style = $s
$ws.Range("E9").Value = '  -5.39%  '
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.112'
$c.Style = $s
$ws.Range("E10").Value = '  -0.87%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '41.05'
$c.Style = $s
$ws.Range("E11").Value = '  -1.96%  '
$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.339'
$c.Style = $s
$ws.Range("E12").Value = '  -0.55%  '
$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = $s
$ws.Range("E13").Value = '  -0.05%  '
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '20.44'
$c.Style = $s
$ws.Range("E14").Value = '  -3.18%  '
$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.811.79'
$c.Style = $s
$ws.Range("E15").Value = '  -0.63%  '
$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.307'
$c.Style = $s
$ws.Range("E16").Value = '  -2.23%  '
$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '92.56'
$c.Style = $s
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("E18").Value = '  -3.24%  '
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.06556'
$c.Style = $s
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -2.74%  '
$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.016'
$c.Style = $s
$ws.Range("E22").Value = '  -0.59%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '28.343.71'
$c.Style = $s
$ws.Range("E23").Value = '  -0.48%  '
$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.13'
$c.Style = $s
$ws.Range("E24").Value = '  -2.82%  '
$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.230'
$c.Style = $s
$ws.Range("E25").Value = '  -0.68%  '
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '161.05'
$c.Style = $s
$ws.Range("E26").Value = '  +1.08%  '
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '20.53'
$c.Style = $s
$ws.Range("E27").Value = '  -2.72%  '
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.021.26'
$c.Style = $s
$ws.Range("E28").Value = '  -0.61%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.418'
$c.Style = $s
$ws.Range("E29").Value = '  +0.63%  '
$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '127.88'
$c.Style = $s
$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.1102'
$c.Style = $s
$ws.Range("E32").Value = '  -2.54%  '
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.671'
$c.Style = $s
$ws.Range("E33").Value = '  -0.38%  '
$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.569'
$c.Style = $s
$ws.Range("E34").Value = '  -2.75%  '
$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.07196'
$c.Style = $s
$ws.Range("E35").Value = '  -4.31%  '
$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.158'
$c.Style = $s
$ws.Range("E36").Value = '  +4.55%  '
$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.02353'
$c.Style = $s
$ws.Range("E37").Value = '  -0.23%  '
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.2187'
$c.Style = $s
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.58'
$c.Style = $s
$ws.Range("E39").Value = '  -4.30%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.050'
$c.Style = $s
$ws.Range("E40").Value = '  -3.92%  '
$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.6184'
$c.Style = $s
$ws.Range("E41").Value = '  -3.20%  '
$ws.Range("E42").Value = '  -0.06%  '
$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.158'
$c.Style = $s
$ws.Range("E43").Value = '  -2.48%  '
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '13.23'
$c.Style = $s
$ws.Range("E44").Value = '  -2.28%  '
$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.5985'
$c.Style = $s
$ws.Range("E45").Value = '  -3.36%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.749'
$c.Style = $s
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.302'
$c.Style = $s
$ws.Range("E47").Value = '  -6.72%  '
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '125.36'
$c.Style = $s
$ws.Range("E48").Value = '  -1.57%  '
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.219'
$c.Style = $s
$ws.Range("E49").Value = '  +1.25%  '
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.921'
$c.Style = $s
$ws.Range("E50").Value = '  -4.07%  '
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.06837'
$c.Style = $s
$ws.Range("E51").Value = '  -1.75%  '
